$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Accept the previously open review points
$ws.Range("E8").Value = "Accepted"
$ws.Range("E10").Value = "Accepted"
$ws.Range("E11").Value = "Accepted"

# Update view: scroll/zoom and move the selection to the last edited cell
$excel.ActiveWindow.Zoom = 83
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("E11").Select()
